# Update column G ("K") values for rows 2-25 with new strikeout counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 7
    4  = 6
    5  = 4
    6  = 4
    7  = 9
    8  = 4
    9  = 8
    10 = 9
    11 = 4
    12 = 6
    13 = 5
    14 = 6
    15 = 9
    16 = 2
    17 = 5
    18 = 9
    19 = 7
    20 = 6
    21 = 9
    22 = 5
    23 = 12
    24 = 5
    25 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
